$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "'003550"
$ws.Range("B8").Value = "LG"
$ws.Range("C8").Value = 79400.0
$ws.Range("D8").Value = 0.0025
$ws.Range("E8").Value = 2.16
$ws.Range("F8").Value = 8.5
$ws.Range("G8").Value = 5.18
$ws.Range("H8").Value = 174674.22
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 33.0
$ws.Range("K8").Value = 33.0

$ws.Rows("14:14").Insert()
$ws.Range("A14").Value = "'005850"
$ws.Range("B14").Value = "에스엘"
$ws.Range("C14").Value = 34150.0
$ws.Range("D14").Value = -0.0353
$ws.Range("E14").Value = 17.42
$ws.Range("F14").Value = 17.22
$ws.Range("G14").Value = 15.37
$ws.Range("H14").Value = 50552.95
$ws.Range("I14").Value = 3.51
$ws.Range("J14").Value = 63.0
$ws.Range("K14").Value = 63.0

$ws.Rows("16:16").Insert()
$ws.Range("A16").Value = "'007340"
$ws.Range("B16").Value = "DN오토모티브"
$ws.Range("C16").Value = 31100.0
$ws.Range("D16").Value = 0.0507
$ws.Range("E16").Value = 19.78
$ws.Range("F16").Value = 17.25
$ws.Range("G16").Value = 23.29
$ws.Range("H16").Value = 31156.34
$ws.Range("I16").Value = 3.22
$ws.Range("J16").Value = 98.0
$ws.Range("K16").Value = 98.0

$ws.Rows("17:17").Insert()
$ws.Range("A17").Value = "'009970"
$ws.Range("B17").Value = "영원무역홀딩스"
$ws.Range("C17").Value = 133900.0
$ws.Range("D17").Value = -0.003
$ws.Range("E17").Value = 13.17
$ws.Range("F17").Value = 10.22
$ws.Range("G17").Value = 18.01
$ws.Range("H17").Value = 201114.72
$ws.Range("I17").Value = 4.0
$ws.Range("J17").Value = 91.0
$ws.Range("K17").Value = 91.0

$ws.Rows("18:18").Insert()
$ws.Range("A18").Value = "'012750"
$ws.Range("B18").Value = "에스원"
$ws.Range("C18").Value = 70900.0
$ws.Range("D18").Value = 0.0028
$ws.Range("E18").Value = 10.85
$ws.Range("F18").Value = 10.27
$ws.Range("G18").Value = 11.06
$ws.Range("H18").Value = 46509.73
$ws.Range("I18").Value = 3.81
$ws.Range("J18").Value = 40.0
$ws.Range("K18").Value = 40.0

$ws.Rows("22:22").Insert()
$ws.Range("A22").Value = "'021240"
$ws.Range("B22").Value = "코웨이"
$ws.Range("C22").Value = 107300.0
$ws.Range("D22").Value = -0.0138
$ws.Range("E22").Value = 19.38
$ws.Range("F22").Value = 17.49
$ws.Range("G22").Value = 20.04
$ws.Range("H22").Value = 43969.66
$ws.Range("I22").Value = 2.45
$ws.Range("J22").Value = 90.0
$ws.Range("K22").Value = 90.0

$ws.Rows("23:23").Insert()
$ws.Range("A23").Value = "'023590"
$ws.Range("B23").Value = "다우기술"
$ws.Range("C23").Value = 36050.0
$ws.Range("D23").Value = 0.03
$ws.Range("E23").Value = 13.16
$ws.Range("F23").Value = 14.03
$ws.Range("G23").Value = 13.29
$ws.Range("H23").Value = 65635.93
$ws.Range("I23").Value = 3.88
$ws.Range("J23").Value = 76.0
$ws.Range("K23").Value = 76.0

$ws.Rows("29:29").Insert()
$ws.Range("A29").Value = "'032830"
$ws.Range("B29").Value = "삼성생명"
$ws.Range("C29").Value = 128500.0
$ws.Range("D29").Value = 0.0023
$ws.Range("E29").Value = 5.76
$ws.Range("F29").Value = 8.44
$ws.Range("G29").Value = 5.33
$ws.Range("H29").Value = 158106.71
$ws.Range("I29").Value = 3.5
$ws.Range("J29").Value = 82.0
$ws.Range("K29").Value = 82.0

$ws.Rows("31:31").Insert()
$ws.Range("A31").Value = "'035250"
$ws.Range("B31").Value = "강원랜드"
$ws.Range("C31").Value = 19480.0
$ws.Range("D31").Value = -0.0041
$ws.Range("E31").Value = 12.08
$ws.Range("F31").Value = 7.85
$ws.Range("G31").Value = 8.39
$ws.Range("H31").Value = 18381.9
$ws.Range("I31").Value = 6.01
$ws.Range("J31").Value = 38.0
$ws.Range("K31").Value = 38.0

$ws.Rows("34:34").Insert()
$ws.Range("A34").Value = "'051600"
$ws.Range("B34").Value = "한전KPS"
$ws.Range("C34").Value = 51300.0
$ws.Range("D34").Value = -0.0039
$ws.Range("E34").Value = 13.26
$ws.Range("F34").Value = 3.52
$ws.Range("G34").Value = 11.65
$ws.Range("H34").Value = 27139.93
$ws.Range("I34").Value = 4.81
$ws.Range("J34").Value = 64.0
$ws.Range("K34").Value = 64.0

$ws.Rows("37:37").Insert()
$ws.Range("A37").Value = "'086280"
$ws.Range("B37").Value = "현대글로비스"
$ws.Range("C37").Value = 147800.0
$ws.Range("D37").Value = 0.0394
$ws.Range("E37").Value = 13.24
$ws.Range("F37").Value = 17.97
$ws.Range("G37").Value = 15.51
$ws.Range("H37").Value = 119232.22
$ws.Range("I37").Value = 2.5
$ws.Range("J37").Value = 97.0
$ws.Range("K37").Value = 97.0

$ws.Rows("40:40").Insert()
$ws.Range("A40").Value = "'111770"
$ws.Range("B40").Value = "영원무역"
$ws.Range("C40").Value = 59200.0
$ws.Range("D40").Value = -0.015
$ws.Range("E40").Value = 12.32
$ws.Range("F40").Value = 9.67
$ws.Range("G40").Value = 18.72
$ws.Range("H40").Value = 84441.54
$ws.Range("I40").Value = 2.36
$ws.Range("J40").Value = 81.0
$ws.Range("K40").Value = 81.0

$ws.Rows("43:43").Insert()
$ws.Range("A43").Value = "'161390"
$ws.Range("B43").Value = "한국타이어앤테크놀로"
$ws.Range("C43").Value = 46100.0
$ws.Range("D43").Value = -0.0254
$ws.Range("E43").Value = 10.79
$ws.Range("F43").Value = 11.5
$ws.Range("G43").Value = 8.93
$ws.Range("H43").Value = 90358.17
$ws.Range("I43").Value = 4.34
$ws.Range("J43").Value = 55.0
$ws.Range("K43").Value = 55.0

$ws.Columns("B:B").ColumnWidth = 20
